# Updated cryptos list on Wed Aug 21 22:46:31 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.974.86"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "2.608.95"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'570.89"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'143.36"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "2.635.68"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("E13").Value = "  +7.05%  "
$ws.Range("D14").Value = "3.074.78"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "60.975.25"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "'23.54"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "2.622.77"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'11.23"
$ws.Range("E19").Value = "  +9.52%  "
$ws.Range("D20").Value = "'4.66"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "'350.33"
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").Value = "'7.11"
$ws.Range("E22").Value = "  +13.20%  "
$ws.Range("D24").Value = "'0.519"
$ws.Range("E24").Value = "  +12.68%  "
$ws.Range("D25").Value = "'64.35"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "'7.71"
$ws.Range("E28").Value = "  +6.11%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("E30").Value = "  +8.35%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.35"
$ws.Range("E32").Value = "  +4.84%  "
$ws.Range("D33").Value = "'160.23"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").Value = "'19.54"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  +5.97%  "
$ws.Range("D36").Value = "'0.973"
$ws.Range("E36").Value = "  +10.85%  "
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = "  +7.34%  "
$ws.Range("D39").Value = "'37.81"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").Value = "'0.855"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("D42").Value = "'299.25"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +7.32%  "
$ws.Range("D44").Value = "'0.0988"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "'0.996"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'0.0548"
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("D48").Value = "'0.0241"
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("D49").Value = "'4.94"
$ws.Range("E49").Value = "  +9.50%  "
$ws.Range("D50").Value = "'10.70"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "'19.68"
$ws.Range("E51").Value = "  +5.93%  "
